$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) cells we touch to remain Text so that
# numeric-looking strings (e.g. "41.93", "1.000", "0.5050") are not
# auto-coerced into numbers by Excel, matching the inline-string cells
# in the original workbook.
$dCells = @("D2","D3","D4","D5","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "25.762.98"
$ws.Range("E2").Value = "  -2.68%  "
$ws.Range("D3").Value = "1.741.94"
$ws.Range("E3").Value = "  -5.12%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "237.76"
$ws.Range("E5").Value = "  -8.83%  "
$ws.Range("D7").Value = "0.5050"
$ws.Range("E7").Value = "  -5.96%  "
$ws.Range("D8").Value = "41.93"
$ws.Range("E8").Value = "  -6.55%  "
$ws.Range("D9").Value = "0.2650"
$ws.Range("E9").Value = "  -12.27%  "
$ws.Range("D10").Value = "0.06157"
$ws.Range("E10").Value = "  -10.37%  "
$ws.Range("D11").Value = "1.743.65"
$ws.Range("E11").Value = "  -5.01%  "
$ws.Range("D12").Value = "0.06921"
$ws.Range("E12").Value = "  -4.28%  "
$ws.Range("D13").Value = "15.34"
$ws.Range("E13").Value = "  -12.87%  "
$ws.Range("D14").Value = "4.491"
$ws.Range("E14").Value = "  -9.68%  "
$ws.Range("D15").Value = "0.5981"
$ws.Range("E15").Value = "  -19.11%  "
$ws.Range("D16").Value = "76.80"
$ws.Range("E16").Value = "  -14.12%  "
$ws.Range("D17").Value = "1.000"
$ws.Range("E17").Value = "  -0.05%  "
$ws.Range("D19").Value = "25.755.58"
$ws.Range("E19").Value = "  -2.81%  "
$ws.Range("D20").Value = "0.000006821"
$ws.Range("E20").Value = "  -13.38%  "
$ws.Range("D21").Value = "11.61"
$ws.Range("E21").Value = "  -16.08%  "
$ws.Range("D22").Value = "1.962.15"
$ws.Range("E22").Value = "  -5.59%  "
$ws.Range("D23").Value = "4.042"
$ws.Range("E23").Value = "  -11.76%  "
$ws.Range("D24").Value = "5.178"
$ws.Range("E24").Value = "  -13.17%  "
$ws.Range("D25").Value = "8.106"
$ws.Range("E25").Value = "  -12.28%  "
$ws.Range("D26").Value = "138.04"
$ws.Range("E26").Value = "  -3.28%  "
$ws.Range("D27").Value = "1.514"
$ws.Range("E27").Value = "  -9.98%  "
$ws.Range("D28").Value = "14.98"
$ws.Range("E28").Value = "  -11.52%  "
$ws.Range("D29").Value = "1.806"
$ws.Range("E29").Value = "  -17.45%  "
$ws.Range("D30").Value = "103.41"
$ws.Range("E30").Value = "  -6.25%  "
$ws.Range("D31").Value = "3.756"
$ws.Range("E31").Value = "  -11.06%  "
$ws.Range("D32").Value = "0.08092"
$ws.Range("E32").Value = "  -8.12%  "
$ws.Range("D33").Value = "3.464"
$ws.Range("E33").Value = "  -13.86%  "
$ws.Range("D34").Value = "0.04515"
$ws.Range("E34").Value = "  -6.24%  "
$ws.Range("D35").Value = "0.9992"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").Value = "2.649"
$ws.Range("E36").Value = "  -9.47%  "
$ws.Range("D37").Value = "0.9782"
$ws.Range("E37").Value = "  -13.47%  "
$ws.Range("E38").Value = "  -16.85%  "
$ws.Range("E39").Value = "  -14.41%  "
$ws.Range("D40").Value = "0.01550"
$ws.Range("E40").Value = "  -9.11%  "
$ws.Range("D41").Value = "1.000"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").Value = "1.899"
$ws.Range("E42").Value = "  -16.90%  "
$ws.Range("D43").Value = "102.93"
$ws.Range("E43").Value = "  -4.59%  "
$ws.Range("D44").Value = "0.3804"
$ws.Range("E44").Value = "  -19.33%  "
$ws.Range("D45").Value = "5.087"
$ws.Range("E45").Value = "  -13.47%  "
$ws.Range("D46").Value = "0.7314"
$ws.Range("E46").Value = "  -19.47%  "
$ws.Range("D47").Value = "0.05346"
$ws.Range("E47").Value = "  -7.58%  "
$ws.Range("D48").Value = "0.1112"
$ws.Range("E48").Value = "  -9.90%  "
$ws.Range("D49").Value = "30.15"
$ws.Range("E49").Value = "  -13.18%  "
$ws.Range("D50").Value = "5.899"
$ws.Range("E50").Value = "  -19.83%  "
$ws.Range("D51").Value = "52.52"

Write-Host "Done applying cryptos update"
